$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Update cell values (corrected / re-typed test data) ---------------
# NOTE: the order cells are first assigned a brand-new text value below is
# deliberate -- it reproduces the order new entries were appended to the
# workbook's shared-string table by the original edit.
$ws.Range("B1").Value = "Payable"
$ws.Range("C1").Value = "Testing1@"

$ws.Range("B4").Value = "TechBite"
$ws.Range("H5").Value = "Workbooks"
$ws.Range("Q4").Value = "Account"

$ws.Range("X9").Value = "201-999-5654"
$ws.Range("Y9").Value = "210-339-0102"
$ws.Range("S9").Value = "Street Sacramento"
$ws.Range("U9").Value = "USA"
$ws.Range("T9").Value = "Alaska"
$ws.Range("V9").Value = 99501

$ws.Range("B5").Value = "TechBite"
$ws.Range("B6").Value = "TechBite"

$ws.Range("S10").Value = "Street Sacramento"
$ws.Range("T10").Value = "Alaska"
$ws.Range("U10").Value = "USA"
$ws.Range("V10").Value = 99501
$ws.Range("X10").Value = "201-999-5654"
$ws.Range("Y10").Value = "210-339-0102"

# --- Update the view state (scroll position + active selection) --------
$win = $excel.ActiveWindow
$win.ScrollColumn = 18
$win.ScrollRow = 1
$ws.Range("W9").Select()
